$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the value in B3
$ws.Range("B3").Value = 54.11

# Update B29's formula to include the IF logic
$ws.Range("B29").Formula = "=IF(B3<100,B17+B25,IF(B3<140,B17+B26,IF(B3<190,B17+B27,B17+B28)))"

# Update sheet view: remove frozen/topLeftCell offset, update zoom, update selection
$excel.ActiveWindow.Zoom = 70
$ws.Range("B52").Select()
